# Add team record (Wins / Losses / Ties) columns to the roster sheet.
# New header cells AD1:AF1 mirror the existing bold/bordered header style
# (copied from AC1), and every data row (2-52, including the repeated
# header row 52) gets the constant team record: 63 wins, 99 losses, 0 ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastDataRow = 52

# --- Header row (row 1): new labels, same formatting as existing headers ---
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows: constant team record on every row ---
for ($r = 2; $r -le $lastDataRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 63   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 99   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
